# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime refreshed with the handback run's timestamp
$wsZhCn.Range("K2").Value = "2016-08-29 21:03:38"
$wsDeDe.Range("K2").Value = "2016-08-29 21:03:45"

# Error Detail cleared now that the handback is in sync
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# Widen the Status columns and shrink the now-empty Error Detail columns
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
